$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the previous question row (row 2) onto the new
# row 4 first, then fill in the values for the new Bible quiz question.
$ws.Range("A2:B2").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A4").Value = "Mark 5:42"
$ws.Range("B4").Value = (Get-Date -Year 2020 -Month 3 -Day 12).Date

# Move the selection to A5, ready for the next question, matching the
# saved workbook state.
$ws.Range("A5").Select()
